$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Add hidden defined name TRNR_... -> #REF! (Bloomberg add-in artifact left by a refresh/upload)
$trnrName = $wb.Names.Add("TRNR_183d342ed17e4d4aaf4f9159f4ea168b_239_1", "#REF!")
$trnrName.Visible = $false

# 2. Replace the USDINR (column F) #N/A N/A placeholders with the refreshed Bloomberg BDH values
$ws.Range("F7").Value = 43.601300000000002
$ws.Range("F8").Value = 43.58
$ws.Range("F9").Value = 43.612000000000002
$ws.Range("F10").Value = 43.645000000000003
$ws.Range("F11").Value = 44.567500000000003
$ws.Range("F12").Value = 44.672499999999999
$ws.Range("F13").Value = 45.14
$ws.Range("F14").Value = 45.75
$ws.Range("F15").Value = 46.037500000000001
$ws.Range("F16").Value = 46.755000000000003
$ws.Range("F17").Value = 46.86
$ws.Range("F18").Value = 46.68
$ws.Range("F19").Value = 46.4
$ws.Range("F20").Value = 46.545000000000002
$ws.Range("F21").Value = 46.64
$ws.Range("F22").Value = 46.83
$ws.Range("F23").Value = 47
$ws.Range("F24").Value = 47.13
$ws.Range("F25").Value = 47.11
$ws.Range("F26").Value = 47.137
$ws.Range("F27").Value = 47.98
$ws.Range("F28").Value = 48
$ws.Range("F29").Value = 47.94
$ws.Range("F30").Value = 48.22
$ws.Range("F31").Value = 48.51
$ws.Range("F32").Value = 48.71
$ws.Range("F33").Value = 48.78
$ws.Range("F34").Value = 48.93
$ws.Range("F35").Value = 49.02
$ws.Range("F36").Value = 48.81
$ws.Range("F37").Value = 48.63
$ws.Range("F38").Value = 48.45
$ws.Range("F39").Value = 48.34
$ws.Range("F40").Value = 48.31
$ws.Range("F41").Value = 48.27
$ws.Range("F42").Value = 47.95
$ws.Range("F43").Value = 47.75
$ws.Range("F44").Value = 47.66
$ws.Range("F45").Value = 47.45
$ws.Range("F46").Value = 47.33
$ws.Range("F47").Value = 47.1
$ws.Range("F48").Value = 46.435000000000002
$ws.Range("F49").Value = 46.16
$ws.Range("F50").Value = 45.83
$ws.Range("F51").Value = 45.55
$ws.Range("F52").Value = 45.31
$ws.Range("F53").Value = 45.695
$ws.Range("F54").Value = 45.625
$ws.Range("F55").Value = 45.28
$ws.Range("F56").Value = 45.2
$ws.Range("F57").Value = 43.35
$ws.Range("F58").Value = 44.83
$ws.Range("F59").Value = 45.41
$ws.Range("F60").Value = 45.85
$ws.Range("F61").Value = 46.305
$ws.Range("F62").Value = 46.305
$ws.Range("F63").Value = 45.83
$ws.Range("F64").Value = 45.375
$ws.Range("F65").Value = 44.27
$ws.Range("F66").Value = 43.414999999999999
$ws.Range("F67").Value = 43.7
$ws.Range("F68").Value = 43.7
$ws.Range("F69").Value = 43.645000000000003
$ws.Range("F70").Value = 43.575000000000003
$ws.Range("F71").Value = 43.75
$ws.Range("F72").Value = 43.465000000000003
$ws.Range("F73").Value = 43.39
$ws.Range("F74").Value = 43.965000000000003
$ws.Range("F75").Value = 44.075000000000003
$ws.Range("F76").Value = 45.075000000000003
$ws.Range("F77").Value = 46.08
$ws.Range("F78").Value = 45.01
$ws.Range("F79").Value = 44.14
$ws.Range("F80").Value = 44.37
$ws.Range("F81").Value = 44.45
$ws.Range("F82").Value = 44.79
$ws.Range("F83").Value = 46.344999999999999
$ws.Range("F84").Value = 46.07
$ws.Range("F85").Value = 46.582500000000003
$ws.Range("F86").Value = 46.52
$ws.Range("F87").Value = 45.924999999999997
$ws.Range("F88").Value = 44.87
$ws.Range("F89").Value = 44.66
$ws.Range("F90").Value = 44.26
$ws.Range("F91").Value = 44.115000000000002
$ws.Range("F92").Value = 44.267000000000003
$ws.Range("F93").Value = 43.465000000000003
$ws.Range("F94").Value = 41.1875
$ws.Range("F95").Value = 40.524999999999999
$ws.Range("F96").Value = 40.655000000000001
$ws.Range("F97").Value = 40.454999999999998
$ws.Range("F98").Value = 40.884999999999998
$ws.Range("F99").Value = 39.85
$ws.Range("F100").Value = 39.3125
$ws.Range("F101").Value = 39.494999999999997
$ws.Range("F102").Value = 39.414999999999999
$ws.Range("F103").Value = 39.354999999999997
$ws.Range("F104").Value = 40.380000000000003
$ws.Range("F105").Value = 40.119999999999997
$ws.Range("F106").Value = 40.484999999999999
$ws.Range("F107").Value = 42.405000000000001
$ws.Range("F108").Value = 43.335000000000001
$ws.Range("F109").Value = 42.354999999999997
$ws.Range("F110").Value = 44.185000000000002
$ws.Range("F111").Value = 46.625
$ws.Range("F112").Value = 48.655000000000001
$ws.Range("F113").Value = 50.274999999999999
$ws.Range("F114").Value = 48.72
$ws.Range("F115").Value = 48.924999999999997
$ws.Range("F116").Value = 51.914999999999999
$ws.Range("F117").Value = 50.734999999999999
$ws.Range("F118").Value = 50.034999999999997
$ws.Range("F119").Value = 46.965000000000003
$ws.Range("F120").Value = 47.89
$ws.Range("F121").Value = 47.634999999999998
$ws.Range("F122").Value = 49.04
$ws.Range("F123").Value = 47.74
$ws.Range("F124").Value = 46.965000000000003
$ws.Range("F125").Value = 46.3125
$ws.Range("F126").Value = 46.534999999999997
$ws.Range("F127").Value = 46.375
$ws.Range("F128").Value = 46.104999999999997
$ws.Range("F129").Value = 44.895000000000003
$ws.Range("F130").Value = 44.52
$ws.Range("F131").Value = 47.155000000000001
$ws.Range("F132").Value = 46.585000000000001
$ws.Range("F133").Value = 46.244999999999997
$ws.Range("F134").Value = 46.81
$ws.Range("F135").Value = 44.475000000000001
$ws.Range("F136").Value = 44.47
$ws.Range("F137").Value = 45.375
$ws.Range("F138").Value = 44.725000000000001
$ws.Range("F139").Value = 45.765000000000001
$ws.Range("F140").Value = 44.945
$ws.Range("F141").Value = 44.594999999999999
$ws.Range("F142").Value = 44.335000000000001
$ws.Range("F143").Value = 44.844999999999999
$ws.Range("F144").Value = 44.582500000000003
$ws.Range("F145").Value = 44.08
$ws.Range("F146").Value = 46.05
$ws.Range("F147").Value = 49.155000000000001
$ws.Range("F148").Value = 49.274999999999999
$ws.Range("F149").Value = 51.465000000000003
$ws.Range("F150").Value = 53.104999999999997
$ws.Range("F151").Value = 49.27
$ws.Range("F152").Value = 49.215000000000003
$ws.Range("F153").Value = 50.945
$ws.Range("F154").Value = 52.7
$ws.Range("F155").Value = 55.89
$ws.Range("F156").Value = 55.594999999999999
$ws.Range("F157").Value = 55.484999999999999
$ws.Range("F158").Value = 55.534999999999997
$ws.Range("F159").Value = 52.625
$ws.Range("F160").Value = 53.76
$ws.Range("F161").Value = 54.674999999999997
$ws.Range("F162").Value = 54.79
$ws.Range("F163").Value = 53.19
$ws.Range("F164").Value = 54.905000000000001
$ws.Range("F165").Value = 54.35
$ws.Range("F166").Value = 53.875
$ws.Range("F167").Value = 56.8
$ws.Range("F168").Value = 59.431199999999997
$ws.Range("F169").Value = 60.265000000000001
$ws.Range("F170").Value = 66.105000000000004
$ws.Range("F171").Value = 62.435000000000002
$ws.Range("F172").Value = 61.844999999999999
$ws.Range("F173").Value = 62.26
$ws.Range("F174").Value = 61.854999999999997
$ws.Range("F175").Value = 62.692500000000003
$ws.Range("F176").Value = 62.03
$ws.Range("F177").Value = 59.722499999999997
$ws.Range("F178").Value = 60.314999999999998
$ws.Range("F179").Value = 59.094999999999999
$ws.Range("F180").Value = 60.085000000000001
$ws.Range("F181").Value = 61.07
$ws.Range("F182").Value = 60.505000000000003
$ws.Range("F183").Value = 61.625
$ws.Range("F184").Value = 61.416200000000003
$ws.Range("F185").Value = 62.078699999999998
$ws.Range("F186").Value = 63.122500000000002
$ws.Range("F187").Value = 61.865000000000002
$ws.Range("F188").Value = 61.952500000000001
$ws.Range("F189").Value = 62.581200000000003
$ws.Range("F190").Value = 63.515000000000001
$ws.Range("F191").Value = 63.69
$ws.Range("F192").Value = 63.615000000000002
$ws.Range("F193").Value = 64.036199999999994
$ws.Range("F194").Value = 66.245000000000005
$ws.Range("F195").Value = 65.564999999999998
$ws.Range("F196").Value = 65.584999999999994
$ws.Range("F197").Value = 66.55
$ws.Range("F198").Value = 66.156300000000002
$ws.Range("F199").Value = 67.844999999999999
$ws.Range("F200").Value = 67.924999999999997
$ws.Range("F201").Value = 66.222499999999997
$ws.Range("F202").Value = 66.415000000000006
$ws.Range("F203").Value = 67.422499999999999
$ws.Range("F204").Value = 67.344999999999999
$ws.Range("F205").Value = 66.765000000000001
$ws.Range("F206").Value = 66.995000000000005
$ws.Range("F207").Value = 66.537499999999994
$ws.Range("F208").Value = 66.737499999999997
$ws.Range("F209").Value = 68.371200000000002
$ws.Range("F210").Value = 67.87
$ws.Range("F211").Value = 67.537499999999994
$ws.Range("F212").Value = 66.825000000000003
$ws.Range("F213").Value = 65.040000000000006
$ws.Range("F214").Value = 64.33
$ws.Range("F215").Value = 64.542500000000004
$ws.Range("F216").Value = 64.858699999999999
$ws.Range("F217").Value = 64.071200000000005
$ws.Range("F218").Value = 64.022300000000001
$ws.Range("F219").Value = 65.319999999999993
$ws.Range("F220").Value = 64.591399999999993
$ws.Range("F221").Value = 64.48
$ws.Range("F222").Value = 63.827500000000001
$ws.Range("F223").Value = 63.847499999999997
$ws.Range("F224").Value = 65.224999999999994
$ws.Range("F225").Value = 65.221699999999998
$ws.Range("F226").Value = 66.739999999999995
$ws.Range("F227").Value = 67.087500000000006
$ws.Range("F228").Value = 68.738699999999994
$ws.Range("F229").Value = 68.456199999999995
$ws.Range("F230").Value = 71.067499999999995
$ws.Range("F231").Value = 72.864999999999995
$ws.Range("F232").Value = 73.537499999999994
$ws.Range("F233").Value = 70.372500000000002
$ws.Range("F234").Value = 69.814999999999998
$ws.Range("F235").Value = 71.287499999999994
$ws.Range("F236").Value = 70.867500000000007
$ws.Range("F237").Value = 69.275000000000006
$ws.Range("F238").Value = 69.553700000000006
$ws.Range("F239").Value = 69.266199999999998
$ws.Range("F240").Value = 68.953699999999998
$ws.Range("F241").Value = 69.058700000000002
$ws.Range("F242").Value = 71.407499999999999
$ws.Range("F243").Value = 71.067499999999995
$ws.Range("F244").Value = 70.818700000000007

# 3. Restore the saved cursor/selection position
$ws.Range("M243").Select()

